$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    '47+41=',
    '28-0=',
    '38+17=',
    '47-44=',
    '29+43=',
    '16+25=',
    '12-7=',
    '98-3=',
    '51+9=',
    '41+2=',
    '57-0=',
    '30+39=',
    '62+32=',
    '28-12=',
    '72-50=',
    '45+49=',
    '47+20=',
    '50+32=',
    '52+11=',
    '91-84=',
    '30-6=',
    '7+67=',
    '33+4=',
    '54+42=',
    '73-64=',
    '16-10=',
    '58-6=',
    '97-85=',
    '6+92=',
    '2+71=',
    '48+16=',
    '75+4=',
    '73-31=',
    '41+7=',
    '51+9=',
    '61-58=',
    '21+53=',
    '11+2=',
    '84-81=',
    '45-41=',
    '25+24=',
    '55+32=',
    '8+5=',
    '78-30=',
    '87-74=',
    '42+52=',
    '27+53=',
    '31+67=',
    '98-6=',
    '76-39=',
    '71-31=',
    '58-21=',
    '75-55=',
    '36-6=',
    '83-53=',
    '33-4=',
    '56-52=',
    '24+32=',
    '55-45=',
    '32+36=',
    '86+1=',
    '27+68=',
    '52-36=',
    '2+67=',
    '40-5=',
    '71-67=',
    '83-65=',
    '54-3=',
    '17+24=',
    '31+59=',
    '88-59=',
    '58-31=',
    '28+64=',
    '67-43=',
    '91-65=',
    '86-48=',
    '88-38=',
    '65+29=',
    '18+62=',
    '45+46=',
    '98-79=',
    '49+9=',
    '34+4=',
    '98-55=',
    '80-18=',
    '71-11=',
    '47+34=',
    '65-13=',
    '4+8=',
    '18+46=',
    '7+24=',
    '17+59=',
    '5+92=',
    '49-31=',
    '55+44=',
    '81-79=',
    '88-22=',
    '85-57=',
    '63-31=',
    '78-65='
)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}
Write-Host "Done. idx=" $idx
